$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I ("Sector"), shifting Sector and
# everything to its right one column to the right.
$ws.Columns("I").Insert()

# Populate the new "Instrument" column.
$ws.Range("I1").Value = "Instrument"
$ws.Range("I2").Value = "Stock"
$ws.Range("I3").Value = "Stock"

# The old "Sector" column (now shifted to J) changes its values from
# "Software" to "Tech".
$ws.Range("J2").Value = "Tech"
$ws.Range("J3").Value = "Tech"

# Update the selected cell to match the authored state.
$ws.Range("J2").Select()
